$d = $word.ActiveDocument

# 1) Capitalize "eclipse" -> "Eclipse" in "using the eclipse IDE"
$d.Content.Find.Execute(
    "using the eclipse IDE", $true, $true, $false, $false, $false,
    $true, 1, $false, "using the Eclipse IDE", 2) | Out-Null

# 2) Remove "I just finished my first year at this university. " before
#    "I coded in Python..."
$d.Content.Find.Execute(
    "(ICS). I just finished my first year at this university. I coded",
    $true, $true, $false, $false, $false,
    $true, 1, $false, "(ICS). I coded", 2) | Out-Null

# 3) Replace the closing sentence about this year's coursework with the
#    expanded description covering the second year (data structures,
#    algorithms course) and this year's upper division classes.
$d.Content.Find.Execute(
    "using IDLE and Eclipse. This year, I am learning C++ and data structures and will take some upper division courses.",
    $true, $true, $false, $false, $false,
    $true, 1, $false,
    "using IDLE and Eclipse. In my second year, I learned about data structures (in C++) and took an algorithms course (BFS, DFS, dynamic programming, etc). This year, I will take some upper division classes.  ",
    2) | Out-Null
